# RunManager.xlsx - "added standup meeting scripts"
#
# On the "Main" run sheet, swap which of two automated test cases is
# flagged to Execute:
#   - row 11, TC61_VerifyStandupMeeting       : Execute True  -> False
#   - row 13, TC81_VerifyUpdateIssueStatus    : Execute False -> True
# and move the sheet's active selection from D11 to D5.
#
# The two "Execute" cells are plain Text-formatted cells holding the
# literal strings "True"/"False" (not native booleans), so a simple
# Range.Value assignment of the word True/False would get auto-coerced
# to a boolean by the engine. Using Range.Copy to swap the two cells'
# full contents (value + style) avoids that and reproduces the exact
# shared-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$d11 = $ws.Range("D11")
$d13 = $ws.Range("D13")
$scratch = $ws.Range("Z100")

$d13.Copy($scratch)
$d11.Copy($d13)
$scratch.Copy($d11)
$scratch.Clear()

# Move the sheet's active selection from D11 to D5.
$ws.Activate()
$ws.Range("D5").Select()
